# Azureskolen-Workshop#2 - bump the cached "today" date shown on the
# Handout Master and Notes Master date placeholders from 24.01.2019 to
# 25.01.2019 (the presentation's datetimeFigureOut fields are fixed/cached
# text, refreshed here the same way PowerPoint refreshes them through the
# Header & Footer object model).

$p = $ppt.ActivePresentation

# Handout Master "Date Placeholder 2" (field id {54C52192-...})
$handoutMaster = $p.HandoutMaster
$handoutMaster.HeadersFooters.DateAndTime.Text = "25.01.2019"

# Notes Master "Date Placeholder 2" (field id {2696BCB2-...})
$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = "25.01.2019"
